$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.701.53'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '2.461.30'
$ws.Range("E3").Value = '  -1.02%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.99'
$ws.Range("E5").Value = '  -1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.32'
$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -1.60%  '

$ws.Range("E9").Value = '  -0.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("E12").Value = '  -1.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.83'
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("D15").Value = '2.919.84'
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").Value = '62.640.54'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").Value = '2.466.14'
$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.62'
$ws.Range("E18").Value = '  -6.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.73'
$ws.Range("E19").Value = '  -2.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.31'
$ws.Range("E20").Value = '  +2.93%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.14'
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '320.94'
$ws.Range("E22").Value = '  -2.86%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.16'
$ws.Range("E24").Value = '  +3.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.78'
$ws.Range("E25").Value = '  -2.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '640.12'
$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("D27").Value = '2.598.17'
$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("D28").Value = '0.0₃0959'
$ws.Range("E28").Value = '  -3.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.41'
$ws.Range("E30").Value = '  -4.85%  '

$ws.Range("E31").Value = '  -2.83%  '

$ws.Range("E32").Value = '  -2.77%  '

$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("E36").Value = '  -2.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.36'
$ws.Range("E37").Value = '  -2.51%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '150.25'
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.364'
$ws.Range("E39").Value = '  -2.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.52'
$ws.Range("E40").Value = '  -1.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  -2.66%  '

$ws.Range("E42").Value = '  -2.30%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  -2.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '153.34'
$ws.Range("E45").Value = '  -1.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.39'
$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.21'
$ws.Range("E48").Value = '  -1.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.604'
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("E50").Value = '  -1.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0902'
$ws.Range("E51").Value = '  -1.94%  '
